$p = $ppt.ActivePresentation
$master = $p.SlideMaster

# Title placeholder: "Click to edit Master title style" -> pig latin
$titleShape = $master.Shapes.Item(1)
$titleTr = $titleShape.TextFrame.TextRange
$titleTr.Paragraphs(1).Text = "ickclay otay edithay astermay itletay estylay"

# Body placeholder: "Click to edit Master text styles" / level texts -> pig latin
$bodyShape = $master.Shapes.Item(2)
$bodyTr = $bodyShape.TextFrame.TextRange
$bodyTr.Paragraphs(1).Text = "ickclay otay edithay astermay exttay esstylay"
$bodyTr.Paragraphs(2).Text = "econdsay evellay"
$bodyTr.Paragraphs(3).Text = "irdthay evellay"
$bodyTr.Paragraphs(4).Text = "ourthfay evellay"
$bodyTr.Paragraphs(5).Text = "ifthfay evellay"
